$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.561.42'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '1.925.17'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.014'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.68%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.012'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4817'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.63%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4049'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08200'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.008'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.78'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.21%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.094'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.898.76'
$ws.Range('E13').Value = '  -3.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.304'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.59'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06870'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.014'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.64'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.012'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').Value = '29.567.75'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.671'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.189'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.20%  '
$ws.Range('D25').Value = '2.151.92'
$ws.Range('E25').Value = '  -1.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.94'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.421'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.03'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.086'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.67%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.73'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.011'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09599'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.71%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.605'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.565'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.385'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06394'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02280'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.199'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5942'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.73'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.012'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.883'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.07%  '
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.502'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.279'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.43'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.07495'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5542'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.78%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.968'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '117.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.433'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.71%  '
